# Updates cryptos list values (prices and volume %) per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting (values like "1.00" or "89.658.91"
# are formatted strings, not numbers) so they round-trip exactly as text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '89.763.65'
$ws.Cells.Item(2, 5).Value = '  +3.00%  '
$ws.Cells.Item(3, 4).Value = '3.189.81'
$ws.Cells.Item(3, 5).Value = '  +1.65%  '
$ws.Cells.Item(4, 5).Value = '  -0.17%  '
$ws.Cells.Item(5, 4).Value = '217.00'
$ws.Cells.Item(5, 5).Value = '  +7.24%  '
$ws.Cells.Item(6, 4).Value = '644.36'
$ws.Cells.Item(6, 5).Value = '  +6.93%  '
$ws.Cells.Item(7, 4).Value = '0.391'
$ws.Cells.Item(7, 5).Value = '  +5.43%  '
$ws.Cells.Item(8, 4).Value = '0.690'
$ws.Cells.Item(8, 5).Value = '  +4.98%  '
$ws.Cells.Item(9, 5).Value = '  +0.00%  '
$ws.Cells.Item(10, 4).Value = '3.188.80'
$ws.Cells.Item(10, 5).Value = '  +1.83%  '
$ws.Cells.Item(11, 4).Value = '0.571'
$ws.Cells.Item(11, 5).Value = '  +8.24%  '
$ws.Cells.Item(12, 5).Value = '  +1.80%  '
$ws.Cells.Item(13, 4).Value = '0.0000255'
$ws.Cells.Item(13, 5).Value = '  +5.88%  '
$ws.Cells.Item(14, 4).Value = '5.38'
$ws.Cells.Item(14, 5).Value = '  +3.55%  '
$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).Value = '33.23'
$ws.Cells.Item(15, 5).Value = '  +4.84%  '
$ws.Cells.Item(16, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(16, 4).Value = '3.771.65'
$ws.Cells.Item(16, 5).Value = '  +1.32%  '
$ws.Cells.Item(17, 4).Value = '89.497.51'
$ws.Cells.Item(17, 5).Value = '  +3.02%  '
$ws.Cells.Item(18, 4).Value = '3.193.45'
$ws.Cells.Item(18, 5).Value = '  +1.24%  '
$ws.Cells.Item(19, 5).Value = '  +14.62%  '
$ws.Cells.Item(20, 5).Value = '  +75.17%  '
$ws.Cells.Item(21, 4).Value = '13.48'
$ws.Cells.Item(21, 5).Value = '  +1.63%  '
$ws.Cells.Item(22, 4).Value = '435.89'
$ws.Cells.Item(22, 5).Value = '  +6.07%  '
$ws.Cells.Item(23, 4).Value = '8.64'
$ws.Cells.Item(23, 5).Value = '  +2.73%  '
$ws.Cells.Item(24, 4).Value = '5.09'
$ws.Cells.Item(24, 5).Value = '  +0.90%  '
$ws.Cells.Item(25, 4).Value = '5.30'
$ws.Cells.Item(25, 5).Value = '  +4.08%  '
$ws.Cells.Item(26, 4).Value = '11.94'
$ws.Cells.Item(26, 5).Value = '  +1.04%  '
$ws.Cells.Item(27, 5).Value = '  +11.46%  '
$ws.Cells.Item(28, 4).Value = '3.357.93'
$ws.Cells.Item(28, 5).Value = '  +1.49%  '
$ws.Cells.Item(29, 5).Value = '  +0.04%  '
$ws.Cells.Item(30, 2).Value = 'Cronos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(30, 4).Value = '0.159'
$ws.Cells.Item(30, 5).Value = '  -1.42%  '
$ws.Cells.Item(31, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(31, 4).Value = '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.27%  '
$ws.Cells.Item(32, 2).Value = 'Bittensor'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(32, 4).Value = '546.34'
$ws.Cells.Item(32, 5).Value = '  +2.53%  '
$ws.Cells.Item(33, 2).Value = 'dogwifhat'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(33, 4).Value = '4.03'
$ws.Cells.Item(33, 5).Value = '  +36.78%  '
$ws.Cells.Item(34, 4).Value = '8.46'
$ws.Cells.Item(34, 5).Value = '  +3.25%  '
$ws.Cells.Item(35, 5).Value = '  +7.27%  '
$ws.Cells.Item(36, 4).Value = '1.92'
$ws.Cells.Item(36, 5).Value = '  +5.43%  '
$ws.Cells.Item(37, 5).Value = '  +1.96%  '
$ws.Cells.Item(38, 4).Value = '22.40'
$ws.Cells.Item(38, 5).Value = '  +3.04%  '
$ws.Cells.Item(39, 4).Value = '22.39'
$ws.Cells.Item(39, 5).Value = '  +2.88%  '
$ws.Cells.Item(40, 5).Value = '  -2.80%  '
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(41, 5).Value = '  +0.14%  '
$ws.Cells.Item(42, 5).Value = '  -0.06%  '
$ws.Cells.Item(43, 4).Value = '1.93'
$ws.Cells.Item(43, 5).Value = '  +3.12%  '
$ws.Cells.Item(44, 4).Value = '0.374'
$ws.Cells.Item(44, 5).Value = '  +1.80%  '
$ws.Cells.Item(45, 4).Value = '146.07'
$ws.Cells.Item(45, 5).Value = '  -0.93%  '
$ws.Cells.Item(46, 4).Value = '173.69'
$ws.Cells.Item(46, 5).Value = '  +1.76%  '
$ws.Cells.Item(47, 4).Value = '43.66'
$ws.Cells.Item(47, 5).Value = '  +1.67%  '
$ws.Cells.Item(48, 4).Value = '0.759'
$ws.Cells.Item(48, 5).Value = '  +10.20%  '
$ws.Cells.Item(49, 5).Value = '  -1.29%  '
$ws.Cells.Item(50, 4).Value = '1.25'
$ws.Cells.Item(50, 5).Value = '  +0.47%  '
$ws.Cells.Item(51, 4).Value = '0.621'
$ws.Cells.Item(51, 5).Value = '  +6.82%  '
